# The author renamed the "alignment" setting row to "indexing" and
# reworded its description: the bowtie index is *generated*, not merely
# "performed". Reflect that on sheet "Foglio1" (already the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "indexing"
$ws.Range("B7").Value = "TRUE if bowtie index needs to be generated, else FALSE"

# The author's last touched/selected cell ended up on B7.
$ws.Range("B7").Select()
